# Scheduled runner: refresh Universalis market-price columns (H-N) on the
# per-job Leve profit sheets. Item/recipe/leve metadata (A-G) is untouched;
# only the price/profit columns move.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1795.3077
$ws.Range("I12").Value = 799
$ws.Range("J12").Value = 1878.3334
$ws.Range("K12").Value = 799
$ws.Range("L12").Value = 1878.3334
$ws.Range("M12").Value = -629
$ws.Range("N12").Value = -2218.3334

$ws.Range("H86").Value = 4793.4614
$ws.Range("I86").Value = 2133
$ws.Range("J86").Value = 5591.6
$ws.Range("K86").Value = 2133
$ws.Range("L86").Value = 5591.6
$ws.Range("M86").Value = -1010
$ws.Range("N86").Value = -7837.6

$ws.Range("H89").Value = 4793.4614
$ws.Range("I89").Value = 2133
$ws.Range("J89").Value = 5591.6
$ws.Range("K89").Value = 10665
$ws.Range("L89").Value = 27958
$ws.Range("M89").Value = -5049
$ws.Range("N89").Value = -39190

$ws.Range("H98").Value = 2629.5
$ws.Range("I98").Value = 2629.5
$ws.Range("K98").Value = 2629.5
$ws.Range("M98").Value = -1131.5

$ws.Range("H116").Value = 10420.675
$ws.Range("I116").Value = 9809.032999999999
$ws.Range("J116").Value = 12255.6
$ws.Range("K116").Value = 9809.032999999999
$ws.Range("L116").Value = 12255.6
$ws.Range("M116").Value = -6367.032999999999
$ws.Range("N116").Value = -19139.6

$ws.Range("H122").Value = 2629.5
$ws.Range("I122").Value = 2629.5
$ws.Range("K122").Value = 7888.5
$ws.Range("M122").Value = -5438.5

$ws.Range("H132").Value = 28040.541
$ws.Range("I132").Value = 30584.633
$ws.Range("J132").Value = 2917.625
$ws.Range("K132").Value = 91753.899
$ws.Range("L132").Value = 8752.875
$ws.Range("M132").Value = -89223.899
$ws.Range("N132").Value = -13812.875

$ws.Range("H137").Value = 1194271.8
$ws.Range("J137").Value = 1821924.5
$ws.Range("L137").Value = 5465773.5
$ws.Range("N137").Value = -5470873.5

$ws.Range("H138").Value = 2912.625
$ws.Range("I138").Value = 2243.5
$ws.Range("J138").Value = 3581.75
$ws.Range("K138").Value = 6730.5
$ws.Range("L138").Value = 10745.25
$ws.Range("M138").Value = -1590.5
$ws.Range("N138").Value = -21025.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6423.125
$ws.Range("I45").Value = 4957.706
$ws.Range("K45").Value = 4957.706
$ws.Range("M45").Value = -4580.706

$ws.Range("H122").Value = 2941.2334
$ws.Range("I122").Value = 1656.7693
$ws.Range("K122").Value = 4970.3079
$ws.Range("M122").Value = -2520.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1833.5385
$ws.Range("I86").Value = 1556.1
$ws.Range("J86").Value = 2758.3333
$ws.Range("K86").Value = 1556.1
$ws.Range("L86").Value = 2758.3333
$ws.Range("M86").Value = -433.0999999999999
$ws.Range("N86").Value = -5004.3333

$ws.Range("H89").Value = 1833.5385
$ws.Range("I89").Value = 1556.1
$ws.Range("J89").Value = 2758.3333
$ws.Range("K89").Value = 7780.5
$ws.Range("L89").Value = 13791.6665
$ws.Range("M89").Value = -2164.5
$ws.Range("N89").Value = -25023.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7512.143
$ws.Range("I31").Value = 2466.818
$ws.Range("J31").Value = 9824.583000000001
$ws.Range("K31").Value = 2466.818
$ws.Range("L31").Value = 9824.583000000001
$ws.Range("M31").Value = -2171.818
$ws.Range("N31").Value = -10414.583

$ws.Range("H34").Value = 7512.143
$ws.Range("I34").Value = 2466.818
$ws.Range("J34").Value = 9824.583000000001
$ws.Range("K34").Value = 2466.818
$ws.Range("L34").Value = 9824.583000000001
$ws.Range("M34").Value = -2264.818
$ws.Range("N34").Value = -10228.583

$ws.Range("H100").Value = 62000
$ws.Range("J100").Value = 62000
$ws.Range("L100").Value = 62000
$ws.Range("N100").Value = -64164

$ws.Range("H105").Value = 883.2857
$ws.Range("I105").Value = 733.36365
$ws.Range("K105").Value = 733.36365
$ws.Range("M105").Value = 1013.63635

$ws.Range("H107").Value = 564.8889
$ws.Range("I107").Value = 545.17645
$ws.Range("K107").Value = 545.17645
$ws.Range("M107").Value = 1374.82355

$ws.Range("H132").Value = 1599.4
$ws.Range("I132").Value = 1516.0416
$ws.Range("K132").Value = 4548.1248
$ws.Range("M132").Value = -2018.1248

$ws.Range("H137").Value = 70709
$ws.Range("I137").Value = 70709
$ws.Range("K137").Value = 70709
$ws.Range("M137").Value = -65609

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6617.9
$ws.Range("I5").Value = 836.4
$ws.Range("J5").Value = 12399.4
$ws.Range("K5").Value = 2509.2
$ws.Range("L5").Value = 37198.2
$ws.Range("M5").Value = -2397.2
$ws.Range("N5").Value = -37422.2

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws.Range("H134").Value = 2264.5557
$ws.Range("I134").Value = 916.4286
$ws.Range("K134").Value = 2749.2858
$ws.Range("M134").Value = 2320.7142

$ws.Range("H135").Value = 6617.9
$ws.Range("I135").Value = 836.4
$ws.Range("J135").Value = 12399.4
$ws.Range("K135").Value = 7527.599999999999
$ws.Range("L135").Value = 111594.6
$ws.Range("M135").Value = -4992.599999999999
$ws.Range("N135").Value = -116664.6

$ws.Range("H139").Value = 1618.3334
$ws.Range("I139").Value = 1376
$ws.Range("K139").Value = 4128
$ws.Range("M139").Value = 1012

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6737.8
$ws.Range("I70").Value = 5922.25
$ws.Range("K70").Value = 5922.25
$ws.Range("M70").Value = -5652.25

$ws.Range("H73").Value = 6737.8
$ws.Range("I73").Value = 5922.25
$ws.Range("K73").Value = 5922.25
$ws.Range("M73").Value = -4986.25

$ws.Range("H80").Value = 9085.714
$ws.Range("J80").Value = 9699.75
$ws.Range("L80").Value = 9699.75
$ws.Range("N80").Value = -11695.75

$ws.Range("H83").Value = 9085.714
$ws.Range("J83").Value = 9699.75
$ws.Range("L83").Value = 48498.75
$ws.Range("N83").Value = -58482.75

$ws.Range("H122").Value = 4278.9375
$ws.Range("I122").Value = 4471.407
$ws.Range("K122").Value = 13414.221
$ws.Range("M122").Value = -10964.221

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5756.154
$ws.Range("I7").Value = 5523
$ws.Range("K7").Value = 5523
$ws.Range("M7").Value = -5411

$ws.Range("H40").Value = 3462.8667
$ws.Range("I40").Value = 3214.4
$ws.Range("K40").Value = 3214.4
$ws.Range("M40").Value = -3078.4

$ws.Range("H122").Value = 4883.4116
$ws.Range("I122").Value = 4699
$ws.Range("J122").Value = 4894.9375
$ws.Range("K122").Value = 14097
$ws.Range("L122").Value = 14684.8125
$ws.Range("M122").Value = -11647
$ws.Range("N122").Value = -19584.8125

$ws.Range("H126").Value = 5756.154
$ws.Range("I126").Value = 5523
$ws.Range("K126").Value = 16569
$ws.Range("M126").Value = -14099

$ws.Range("H132").Value = 845917.3
$ws.Range("I132").Value = 1049921.8
$ws.Range("J132").Value = 4398.75
$ws.Range("K132").Value = 3149765.4
$ws.Range("L132").Value = 13196.25
$ws.Range("M132").Value = -3147235.4
$ws.Range("N132").Value = -18256.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2704.3333
$ws.Range("I107").Value = 1235.7273
$ws.Range("J107").Value = 6743
$ws.Range("K107").Value = 3707.1819
$ws.Range("L107").Value = 20229
$ws.Range("M107").Value = -1787.1819
$ws.Range("N107").Value = -24069

$ws.Range("H122").Value = 2762.1562
$ws.Range("I122").Value = 2049.3333
$ws.Range("K122").Value = 6147.999899999999
$ws.Range("M122").Value = -3697.999899999999

$ws.Range("H132").Value = 1237532.8
$ws.Range("I132").Value = 1820616.9
$ws.Range("K132").Value = 5461850.699999999
$ws.Range("M132").Value = -5459320.699999999
